# Auto-generated edit script applying the cryptos.xlsx price/listing update
# (commit: "Updated symbol list on Wed Dec 28 16:34:38 UTC 2022 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values but must stay TEXT (t="inlineStr" in the
# source file), matching the original cell type. Force text storage via NumberFormat
# "@" before writing, then restore the default "Normal" style so no stray style index
# is left behind on the cell.
$priceCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9",
    "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17",
    "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25",
    "D26", "D28", "D40", "D41", "D42", "D43", "D44", "D45",
    "D46", "D47", "D48", "D49", "D50"
)
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# --- Row 2-9: price refresh only ---
$ws.Range("D2").Value = "243.69"
$ws.Range("D3").Value = "23.81"
$ws.Range("D4").Value = "5.264"
$ws.Range("D5").Value = "0.05834"
$ws.Range("D6").Value = "6.465"
$ws.Range("D7").Value = "3.334"
$ws.Range("D8").Value = "0.8088"
$ws.Range("D9").Value = "0.8883"

# --- Rows 10-18: exchange-token list shifted down by one (cyclic; "One" moves to
# the top of the block), plus refreshed price/volume-label text ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01034"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1379"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07156"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03086"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03030"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09336"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.815"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001543"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04699"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- Remaining rows: price refresh (+ a couple of "Best/Worst in 24h" label tweaks) ---
$ws.Range("D19").Value = "0.006222"
$ws.Range("D20").Value = "0.001261"
$ws.Range("D21").Value = "0.003851"
$ws.Range("D22").Value = "0.00008701"
$ws.Range("D23").Value = "3.563"
$ws.Range("D24").Value = "2.170"
$ws.Range("D25").Value = "0.3192"
$ws.Range("D26").Value = "0.1305"
$ws.Range("D28").Value = "0.0002340"
$ws.Range("D40").Value = "0.03784"
$ws.Range("D41").Value = "0.006341"
$ws.Range("D42").Value = "0.1051"
$ws.Range("D43").Value = "0.002488"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "0.006934"
$ws.Range("D45").Value = "0.00005314"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.5112"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "0.002164"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").Value = "0.0002000"

# Clear the explicit Text-number-format style back to the workbook default so the
# saved cell keeps style index 0, matching the original.
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}
